$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before DW (column 127), shifting DW:FA -> DX:FB.
$ws.Range("DW1").EntireColumn.Insert()

# Populate the newly inserted column: header "18-nov" in row 1,
# and "-" placeholders for every data row (2-25), matching the
# formatting/style already used by the rest of the sheet.
$ws.Range("DW1").Value = "18-nov"
$ws.Range("DW2:DW25").Value = "-"
